$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the formatting of the previous status row (row 23) onto the new
# row 24, but skip column C (which has no cell in row 23) so we don't
# introduce a stray blank cell.
$ws.Range("A23:B23").Copy()
$ws.Range("A24:B24").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D23").Copy()
$ws.Range("D24").PasteSpecial(-4122)      # xlPasteFormats

# New daily status entry for 2/2/2021 (Excel serial date 44229)
$ws.Range("A24").Value = 44229

$ws.Range("B24").Value = "1. I have completed 4 string programs and 6 recursion programs from the given list of C questions and pushed the programs into git https://github.com/gandepallipavani/C_Programs`n2. I worked on testcases writing for whatsapp as per the task given by Srinivas`n3. I worked on writing sample bug ticket in notepad as per the task given by Srinivas`n4. Attended the session by Srinivas about validating the testcases written`n5. Completed 2 hacker rank programs "

$ws.Range("D24").Value = "PointersInC.txt`nComditionalStatementsInC.txt"

$ws.Rows.Item(24).RowHeight = 180

# Update the saved selection to the newly added row, matching the
# author's last cursor position.
$ws.Range("B24").Select()
